$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.090240666666667
$ws.Range("H2").Value = 15.270722
$ws.Range("I2").Value = 0.01518526656315472
$ws.Range("J2").Value = 0.01525191836740238
$ws.Range("M2").Value = 1.009860666666667
$ws.Range("N2").Value = 3.029582
$ws.Range("O2").Value = 0.01353413605720072
$ws.Range("P2").Value = 0.01542521070970148
$ws.Range("Q2").Value = 5.140433833133779
$ws.Range("R2").Value = 46.26390449820401
$ws.Range("S2").Value = 0.0002055194637305968
$ws.Range("T2").Value = 0.0002352640545443479
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.090240666666667
$ws.Range("H3").Value = 15.270722
$ws.Range("I3").Value = 0.01518526656315472
$ws.Range("J3").Value = 0.01525191836740238
$ws.Range("O3").Value = 0.6185519418990597
$ws.Range("P3").Value = 0.704979911415303
$ws.Range("Q3").Value = 234.9337494650671
$ws.Range("R3").Value = 2114.403745185604
$ws.Range("S3").Value = 0.009392876120894213
$ws.Range("T3").Value = 0.01075229605956476
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.090240666666667
$ws.Range("H4").Value = 15.270722
$ws.Range("I4").Value = 0.01518526656315472
$ws.Range("J4").Value = 0.01525191836740238
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.009315666666666667
$ws.Range("N4").Value = 0.027947
$ws.Range("O4").Value = 0.0001248484115599408
$ws.Range("P4").Value = 0.000142293017222847
$ws.Range("Q4").Value = 0.04741898530377778
$ws.Range("R4").Value = 0.426770867734
$ws.Range("S4").Value = 0.000001895856409524148
$ws.Range("T4").Value = 0.000002170241482934244
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.090240666666667
$ws.Range("H5").Value = 15.270722
$ws.Range("I5").Value = 0.01518526656315472
$ws.Range("J5").Value = 0.01525191836740238
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 27.4428835
$ws.Range("N5").Value = 54.885767
$ws.Range("O5").Value = 0.3677890736321797
$ws.Range("P5").Value = 0.2794525848577725
$ws.Range("Q5").Value = 139.6908816022957
$ws.Range("R5").Value = 838.145289613774
$ws.Range("S5").Value = 0.005584975122120387
$ws.Range("T5").Value = 0.004262188011810334
$ws.Range("I6").Value = 0.003043737298347591
$ws.Range("J6").Value = 0.003057096996825524
$ws.Range("M6").Value = 1.009860666666667
$ws.Range("N6").Value = 3.029582
$ws.Range("O6").Value = 0.01353413605720072
$ws.Range("P6").Value = 0.01542521070970148
$ws.Range("Q6").Value = 1.030349393112445
$ws.Range("R6").Value = 9.273144538012
$ws.Range("S6").Value = 0.00004119435471821285
$ws.Range("T6").Value = 0.0000471563653360293
$ws.Range("I7").Value = 0.003043737298347591
$ws.Range("J7").Value = 0.003057096996825524
$ws.Range("O7").Value = 0.6185519418990597
$ws.Range("P7").Value = 0.704979911415303
$ws.Range("S7").Value = 0.0018827096165235
$ws.Range("T7").Value = 0.002155191970010047
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.003043737298347591
$ws.Range("J8").Value = 0.003057096996825524
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.009315666666666667
$ws.Range("N8").Value = 0.027947
$ws.Range("O8").Value = 0.0001248484115599408
$ws.Range("P8").Value = 0.000142293017222847
$ws.Range("Q8").Value = 0.009504669122444444
$ws.Range("R8").Value = 0.085542022102
$ws.Range("S8").Value = 0.0000003800057669044424
$ws.Range("T8").Value = 0.0000004350035556212081
$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = 0.003043737298347591
$ws.Range("J9").Value = 0.003057096996825524
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 27.4428835
$ws.Range("N9").Value = 54.885767
$ws.Range("O9").Value = 0.3677890736321797
$ws.Range("P9").Value = 0.2794525848577725
$ws.Range("Q9").Value = 27.99966301570367
$ws.Range("R9").Value = 167.997978094222
$ws.Range("S9").Value = 0.001119453321338974
$ws.Range("T9").Value = 0.0008543136579238262
$ws.Range("G10").Value = 69.13821133333333
$ws.Range("H10").Value = 207.414634
$ws.Range("I10").Value = 0.206253935235621
$ws.Range("J10").Value = 0.2071592335956769
$ws.Range("M10").Value = 1.009860666666667
$ws.Range("N10").Value = 3.029582
$ws.Range("O10").Value = 0.01353413605720072
$ws.Range("P10").Value = 0.01542521070970148
$ws.Range("Q10").Value = 69.81996018922089
$ws.Range("R10").Value = 628.3796417029879
$ws.Range("S10").Value = 0.002791468821811961
$ws.Range("T10").Value = 0.003195474828673586
$ws.Range("G11").Value = 69.13821133333333
$ws.Range("H11").Value = 207.414634
$ws.Range("I11").Value = 0.206253935235621
$ws.Range("J11").Value = 0.2071592335956769
$ws.Range("O11").Value = 0.6185519418990597
$ws.Range("P11").Value = 0.704979911415303
$ws.Range("Q11").Value = 3190.988458800087
$ws.Range("R11").Value = 28718.89612920078
$ws.Range("S11").Value = 0.1275787721643163
$ws.Range("T11").Value = 0.1460430981491424
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = 69.13821133333333
$ws.Range("H12").Value = 207.414634
$ws.Range("I12").Value = 0.206253935235621
$ws.Range("J12").Value = 0.2071592335956769
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.009315666666666667
$ws.Range("N12").Value = 0.027947
$ws.Range("O12").Value = 0.0001248484115599408
$ws.Range("P12").Value = 0.000142293017222847
$ws.Range("Q12").Value = 0.6440685307108889
$ws.Range("R12").Value = 5.796616776397999
$ws.Range("S12").Value = 0.00002575047619215419
$ws.Range("T12").Value = 0.00002947731239390144
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 69.13821133333333
$ws.Range("H13").Value = 207.414634
$ws.Range("I13").Value = 0.206253935235621
$ws.Range("J13").Value = 0.2071592335956769
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 27.4428835
$ws.Range("N13").Value = 54.885767
$ws.Range("O13").Value = 0.3677890736321797
$ws.Range("P13").Value = 0.2794525848577725
$ws.Range("Q13").Value = 1897.351879019046
$ws.Range("R13").Value = 11384.11127411428
$ws.Range("S13").Value = 0.07585794377330064
$ws.Range("T13").Value = 0.05789118330546702
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 4.394653
$ws.Range("H14").Value = 8.789306
$ws.Range("I14").Value = 0.01311018115402158
$ws.Range("J14").Value = 0.008778483271329277
$ws.Range("M14").Value = 1.009860666666667
$ws.Range("N14").Value = 3.029582
$ws.Range("O14").Value = 0.01353413605720072
$ws.Range("P14").Value = 0.01542521070970148
$ws.Range("Q14").Value = 4.437987208348667
$ws.Range("R14").Value = 26.627923250092
$ws.Range("S14").Value = 0.0001774349754730768
$ws.Range("T14").Value = 0.0001354099541718436
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 4.394653
$ws.Range("H15").Value = 8.789306
$ws.Range("I15").Value = 0.01311018115402158
$ws.Range("J15").Value = 0.008778483271329277
$ws.Range("O15").Value = 0.6185519418990597
$ws.Range("P15").Value = 0.704979911415303
$ws.Range("Q15").Value = 202.8297627750486
$ws.Range("R15").Value = 1216.978576650292
$ws.Range("S15").Value = 0.008109328011468504
$ws.Range("T15").Value = 0.006188654358982433
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 4.394653
$ws.Range("H16").Value = 8.789306
$ws.Range("I16").Value = 0.01311018115402158
$ws.Range("J16").Value = 0.008778483271329277
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.009315666666666667
$ws.Range("N16").Value = 0.027947
$ws.Range("O16").Value = 0.0001248484115599408
$ws.Range("P16").Value = 0.000142293017222847
$ws.Range("Q16").Value = 0.04093912246366666
$ws.Range("R16").Value = 0.245634734782
$ws.Range("S16").Value = 0.000001636785292342666
$ws.Range("T16").Value = 0.000001249116871317731
$ws.Range("D17").Value = "MuSCs"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 4.394653
$ws.Range("H17").Value = 8.789306
$ws.Range("I17").Value = 0.01311018115402158
$ws.Range("J17").Value = 0.008778483271329277
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 27.4428835
$ws.Range("N17").Value = 54.885767
$ws.Range("O17").Value = 0.3677890736321797
$ws.Range("P17").Value = 0.2794525848577725
$ws.Range("Q17").Value = 120.6019503019255
$ws.Range("R17").Value = 482.407801207702
$ws.Range("S17").Value = 0.004821781381787657
$ws.Range("T17").Value = 0.002453169841303681
$ws.Range("G18").Value = 255.5657806666667
$ws.Range("H18").Value = 766.6973419999999
$ws.Range("I18").Value = 0.7624068797488551
$ws.Range("J18").Value = 0.765753267768766
$ws.Range("M18").Value = 1.009860666666667
$ws.Range("N18").Value = 3.029582
$ws.Range("O18").Value = 0.01353413605720072
$ws.Range("P18").Value = 0.01542521070970148
$ws.Range("Q18").Value = 258.0858296412271
$ws.Range("R18").Value = 2322.772466771044
$ws.Range("S18").Value = 0.01031851844146688
$ws.Range("T18").Value = 0.01181190550697567
$ws.Range("G19").Value = 255.5657806666667
$ws.Range("H19").Value = 766.6973419999999
$ws.Range("I19").Value = 0.7624068797488551
$ws.Range("J19").Value = 0.765753267768766
$ws.Range("O19").Value = 0.6185519418990597
$ws.Range("P19").Value = 0.704979911415303
$ws.Range("Q19").Value = 11795.32187547916
$ws.Range("R19").Value = 106157.8968793124
$ws.Range("S19").Value = 0.4715882559858572
$ws.Range("T19").Value = 0.5398406708776035
$ws.Range("D20").Value = "Inflammatory-Mac"
$ws.Range("G20").Value = 255.5657806666667
$ws.Range("H20").Value = 766.6973419999999
$ws.Range("I20").Value = 0.7624068797488551
$ws.Range("J20").Value = 0.765753267768766
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.009315666666666667
$ws.Range("N20").Value = 0.027947
$ws.Range("O20").Value = 0.0001248484115599408
$ws.Range("P20").Value = 0.000142293017222847
$ws.Range("Q20").Value = 2.380765624097111
$ws.Range("R20").Value = 21.426890616874
$ws.Range("S20").Value = 0.00009518528789901536
$ws.Range("T20").Value = 0.0001089613429190724
$ws.Range("D21").Value = "MuSCs"
$ws.Range("G21").Value = 255.5657806666667
$ws.Range("H21").Value = 766.6973419999999
$ws.Range("I21").Value = 0.7624068797488551
$ws.Range("J21").Value = 0.765753267768766
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 27.4428835
$ws.Range("N21").Value = 54.885767
$ws.Range("O21").Value = 0.3677890736321797
$ws.Range("P21").Value = 0.2794525848577725
$ws.Range("Q21").Value = 7013.461945421886
$ws.Range("R21").Value = 42080.77167253131
$ws.Range("S21").Value = 0.280404920033632
$ws.Range("T21").Value = 0.2139917300412677
